$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Checkout payments sheet: it was the tab-selected sheet; its selection
#    moves to I1 and it stops being tab-selected (Minicart takes over below)
# ---------------------------------------------------------------------------
$wsCheckout = $wb.Worksheets.Item("Checkout payments")
$wsCheckout.Activate()
$wsCheckout.Range("I1").Select()

# ---------------------------------------------------------------------------
# 2. Minicart sheet: restructure data and become the tab-selected sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Minicart")

# Insert two new columns before column F (shifts F:V -> H:X)
$ws.Range("F1:G1").EntireColumn.Insert()

# New header cells for the inserted columns
$ws.Range("F1").Value = "outdoor"
$ws.Range("G1").Value = "Trail"

# New data cells on existing rows
$ws.Range("C3").Value = "$"
$ws.Range("C5").Value = "$"

# Value change on the (post-shift) product-name cell in row 3
$ws.Range("S3").Value = "Heritage Waist Pack 8 NanoFly®"

# New row 6
$ws.Range("A6").Value = "Backpacks & Bags"
$ws.Range("F6").Value = "Outdoor Packs"
$ws.Range("G6").Value = "Hiking Packs"

# I6 / J6: carry the "quote-prefix" text format used elsewhere in the sheet
# (e.g. N2, the post-shift equivalent of the original L2) but stay empty -
# copy format only, then clear the content.
$ws.Range("N2").Copy()
$ws.Range("I6:J6").PasteSpecial(-4122)
$ws.Range("I6:J6").ClearContents()
$excel.CutCopyMode = $false

# Make Minicart the active/tab-selected sheet, with the new selection
$ws.Activate()
$ws.Range("E11").Select()
